$wb = $excel.ActiveWorkbook

function Update-Row($Sheet, $Row) {
    # Columns that already contain non-numeric-looking text (currency symbol) can be
    # assigned directly without risking automatic number coercion.
    $Sheet.Range("D$Row").Value = "¥1,000,689.80"
    $Sheet.Range("E$Row").Value = "¥+689.80"

    # Columns whose text looks like a plain number/percentage must be forced to stay
    # text (leading apostrophe = "treat as text", mirrors typing it in the Excel UI).
    $Sheet.Range("F$Row").Value = "'+0.07%"
    $Sheet.Range("G$Row").Value = "'+1.35%"

    # Numeric cell - stays a real number.
    $Sheet.Range("H$Row").Value = -1.201

    $Sheet.Range("J$Row").Value = "'58.3%"
    $Sheet.Range("K$Row").Value = "'0.0058%"
    $Sheet.Range("L$Row").Value = "'0.4422%"

    # Numeric cells.
    $Sheet.Range("M$Row").Value = 13
    $Sheet.Range("N$Row").Value = 13

    # Text that looks like a pure integer must also be forced to stay text.
    $Sheet.Range("P$Row").Value = "'20260106"
}

$summarySheet = $wb.Worksheets.Item("Summary")
Update-Row $summarySheet 15

$patternSheet = $wb.Worksheets.Item("Pattern3-Data+News")
Update-Row $patternSheet 5
